$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.795.05"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "1.760.18"
$ws.Range("E3").Value = "  -2.96%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'323.26"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.4268"
$ws.Range("E7").Value = "  -3.87%  "
$ws.Range("D8").Value = "'0.3623"
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("D9").Value = "'0.07578"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("D10").Value = "'42.69"
$ws.Range("E10").Value = "  -4.79%  "
$ws.Range("D11").Value = "'1.095"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "'20.67"
$ws.Range("E13").Value = "  -6.13%  "
$ws.Range("D14").Value = "'6.057"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("D15").Value = "'7.274"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").Value = "1.752.03"
$ws.Range("E16").Value = "  -3.67%  "
$ws.Range("D17").Value = "'91.30"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "'0.00001073"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "'0.06380"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'17.07"
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("D22").Value = "'5.916"
$ws.Range("E22").Value = "  -5.01%  "
$ws.Range("D23").Value = "27.856.40"
$ws.Range("E23").Value = "  -1.51%  "
$ws.Range("D24").Value = "'11.21"
$ws.Range("E24").Value = "  -4.17%  "
$ws.Range("D25").Value = "'2.122"
$ws.Range("E25").Value = "  +4.78%  "
$ws.Range("D26").Value = "'160.89"
$ws.Range("E26").Value = "  +3.72%  "
$ws.Range("D27").Value = "'20.34"
$ws.Range("D28").Value = "1.952.57"
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("D29").Value = "'2.136"
$ws.Range("E29").Value = "  -7.95%  "
$ws.Range("D30").Value = "'124.92"
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("D31").Value = "'1.117"
$ws.Range("E31").Value = "  -6.78%  "
$ws.Range("D32").Value = "'3.682"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "'5.584"
$ws.Range("E33").Value = "  -4.70%  "
$ws.Range("D34").Value = "'0.08926"
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("D35").Value = "'12.23"
$ws.Range("E35").Value = "  -6.39%  "
$ws.Range("D36").Value = "'0.02301"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("D37").Value = "'0.2107"
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("D38").Value = "'0.06029"
$ws.Range("E38").Value = "  -2.74%  "
$ws.Range("D39").Value = "'0.6368"
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("D40").Value = "'4.990"
$ws.Range("E40").Value = "  -3.67%  "
$ws.Range("D41").Value = "'1.181"
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.401"
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'7.872"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").Value = "'13.33"
$ws.Range("E45").Value = "  -4.20%  "
$ws.Range("D46").Value = "'0.5888"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").Value = "'3.701"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").Value = "'1.987"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("D49").Value = "'122.70"
$ws.Range("E49").Value = "  -3.05%  "
$ws.Range("D50").Value = "'1.181"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("E51").Value = "  -2.04%  "

# Reset the style on cells that needed a text quote-prefix (to avoid
# Excel auto-converting numeric-looking strings to numbers), restoring
# them to the default/unstyled state matching the target workbook.
# Union() batches are kept <=8 args and row-contiguous (engine constraint).
$excel.Union($ws.Range("D4"), $ws.Range("D5"), $ws.Range("D6"), $ws.Range("D7"), $ws.Range("D8"), $ws.Range("D9"), $ws.Range("D10"), $ws.Range("D11")).Style = "Normal"
$excel.Union($ws.Range("D13"), $ws.Range("D14"), $ws.Range("D15")).Style = "Normal"
$excel.Union($ws.Range("D17"), $ws.Range("D18"), $ws.Range("D19"), $ws.Range("D20"), $ws.Range("D21"), $ws.Range("D22")).Style = "Normal"
$excel.Union($ws.Range("D24"), $ws.Range("D25"), $ws.Range("D26"), $ws.Range("D27")).Style = "Normal"
$excel.Union($ws.Range("D29"), $ws.Range("D30"), $ws.Range("D31"), $ws.Range("D32"), $ws.Range("D33"), $ws.Range("D34"), $ws.Range("D35"), $ws.Range("D36")).Style = "Normal"
$excel.Union($ws.Range("D37"), $ws.Range("D38"), $ws.Range("D39"), $ws.Range("D40"), $ws.Range("D41"), $ws.Range("D42"), $ws.Range("D43"), $ws.Range("D44")).Style = "Normal"
$excel.Union($ws.Range("D45"), $ws.Range("D46"), $ws.Range("D47"), $ws.Range("D48"), $ws.Range("D49"), $ws.Range("D50")).Style = "Normal"
